# Insert a new data row at row 14 ("Fruta / hortaliza, semanal" commit: a new
# weekly price observation was added to the Tuna / Vega Modelo de Temuco
# sheet). Excel shifts every existing row from 14 downward by one, so the
# former row 14 becomes row 15, ..., former row 114 becomes row 115 — this
# matches the dimension growing from A1:T114 to A1:T115.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 14..114 down to 15..115, leaving a blank row 14 to fill in.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new observation.
$ws.Range("A14").Value2 = 10
$ws.Range("B14").Value2 = "Vega Modelo de Temuco"
$ws.Range("C14").Value2 = "La Araucanía"
$ws.Range("D14").Value2 = 45163
$ws.Range("E14").Value2 = 9
$ws.Range("F14").Value2 = "Fruta"
$ws.Range("G14").Value2 = 100107
$ws.Range("H14").Value2 = "Otros"
$ws.Range("I14").Value2 = 100107011
$ws.Range("J14").Value2 = "Tuna"
$ws.Range("K14").Value2 = "Sin especificar"
$ws.Range("L14").Value2 = "Primera"
$ws.Range("M14").Value2 = 50
$ws.Range("N14").Value2 = 32000
$ws.Range("O14").Value2 = 32000
$ws.Range("P14").Value2 = 32000
$ws.Range("Q14").Value2 = "$/caja 16 kilos"
$ws.Range("R14").Value2 = "Provincia de Los Andes"
$ws.Range("S14").Value2 = 2000
$ws.Range("T14").Value2 = 16
